$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two extra contact rows (Virat Kohli, Rohit Sharma), keeping only
# the header row and a single contact row.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Replace the remaining contact (Sam Curran) with the new contact added
# through the API.
$ws.Range("B2").Value = "def"
$ws.Range("C2").Value = "def"
$ws.Range("A2").Value = "def@gmail.com"

# Rebuild the hyperlink collection so only the new contact's e-mail link
# remains (drops the now-orphaned links that used to sit on rows 3 & 4).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:def@gmail.com")
$ws.Range("A2").Style = "Hyperlink"

$ws.Range("C7").Select()
